$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns B, C, D, E, G (F is unchanged)
$data = @{
    2  = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987;  G = 2.964545797025059 }
    3  = @{ B = 0.2881169905109251;  C = 1.626987699542094;  D = 3.223369029078222;   E = 13.86384647080068;   G = 19.00232018993193 }
    4  = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987;  G = 5.582307763322248 }
    5  = @{ B = 1.445647641019636;   C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987;  G = 6.82939032824165 }
    6  = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987;  G = 5.582307763322248 }
    7  = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 3.223369029078222;   E = 13.86384647080068;   G = 19.36876847130326 }
    8  = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987;  G = 3.536033448013082 }
    9  = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987;  G = 6.15379541431027 }
    10 = @{ B = 1.445647641019636;   C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987;  G = 6.82939032824165 }
    11 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987;  G = 6.15379541431027 }
    12 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987;  G = 8.656069925401464 }
    13 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 18.71679738969934;   E = 13.86384647080068;   G = 37.47995879822157 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
